$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (22-30), following the same pattern as the
# existing rows: A = regcntr_id, B = machine_id, C = "eng", D = TRUE,
# E = "superadmin", F/G = "now()"
$startRow = 22
for ($i = 0; $i -lt 9; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10002 + $i
    $ws.Cells.Item($row, 2).Value = 10021 + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Scroll the view down and select the newly added machine_id column
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22:B30").Select() | Out-Null

# Set up the page for printing (portrait)
$ws.PageSetup.Orientation = 1
